$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price report was inserted as row 111 ("Ajo", "Chino", "Primera"
# for the week with date serial 44468). Every row that used to occupy 111..148
# shifts down by one (to 112..149); inserting the row this way reproduces that
# shift automatically.
$ws.Rows.Item(111).Insert()

$ws.Range("A111").Value = 8
$ws.Range("B111").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C111").Value = 'Coquimbo'
$ws.Range("D111").Value = 44468
$ws.Range("E111").Value = 4
$ws.Range("F111").Value = 100112003
$ws.Range("G111").Value = 'Ajo'
$ws.Range("H111").Value = 'Chino'
$ws.Range("I111").Value = 'Primera'
$ws.Range("J111").Value = 660
$ws.Range("K111").Value = 16000
$ws.Range("L111").Value = 17000
$ws.Range("M111").Value = 16500
$ws.Range("N111").Value = '$/caja 10 kilos'
$ws.Range("O111").Value = 'China'
$ws.Range("P111").Value = 1650
$ws.Range("Q111").Value = 10
$ws.Range("R111").Value = 'Hortaliza'
